$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LivingWill")

$ws.Range("K5").Value = "Modeling question: leave it open? e.g. Source[X]"
$ws.Range("K8").Value = "Be more open with regard to references (Patient/Practitioner/Organization)? "
$ws.Range("K6").Value = "No note or comment element available - gForge #13313"

$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H30").Select()
